# Update the keyword co-occurrence cluster table (Tables(1)) with the
# redone clustering results. Columns: 1=Cluster, 2=Key Terms, 3=Size,
# 4=Centrality, 5=Density. Row 1 is the header row, data rows are 2..6.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-Cell($table, $row, $col, $text) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $text
}

# Row 2 -> Cluster 1 (red, D62728)
Set-Cell $t 2 2 "cardiovascular, p300, stress, adult, amplitude, child, performance, adolescent, depression, reduce, age"
Set-Cell $t 2 3 "11"
Set-Cell $t 2 4 "1809 (3)"
Set-Cell $t 2 5 "1994 (2)"

# Row 3 -> Cluster 2 (green, 2CA02C)
Set-Cell $t 3 2 "cognition, error, function, increase, cardiac, time, cortex, mechanism, fear, sensitivity, control"
Set-Cell $t 3 3 "11"
Set-Cell $t 3 4 "2003 (1)"
Set-Cell $t 3 5 "2134 (1)"

# Row 4 -> Cluster 3 (blue, 1F77B4)
Set-Cell $t 4 2 "erp, visual, eeg, human, cue, component, sound, activation, inhibition, paradigm"
Set-Cell $t 4 3 "10"
Set-Cell $t 4 4 "1882 (2)"
Set-Cell $t 4 5 "1800 (4)"

# Row 5 -> Cluster 4 (yellow, BCBD22)
Set-Cell $t 5 2 "attention, electrophysiological, auditory, neural, memory, behavior, mmn, scene, detection, context"
Set-Cell $t 5 3 "10"
Set-Cell $t 5 4 "1669 (5)"
Set-Cell $t 5 5 "1575 (5)"

# Row 6 -> Cluster 5 (purple, 9467BD)
Set-Cell $t 6 2 "emotion, perception, startle, affect, brain, healthy, impact, prepulse inhibition"
Set-Cell $t 6 3 "8"
Set-Cell $t 6 4 "1672 (4)"
Set-Cell $t 6 5 "1800 (3)"
